$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nbsp = [char]0x00A0

# Delete row 11 (the "Elwanda" child row) so every following row shifts up by
# one. This turns the old 14-row sheet into a 13-row sheet (old rows 12-14
# become the new rows 11-13), matching the new dimension A1:H13.
$ws.Rows.Item(11).Delete()

# --- nChildren count dropped from 6 to 5 ---
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "5"

# --- Row 6 (child 0) ---
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "12"
$ws.Range("C6").Value = "Frankie " + $nbsp
$ws.Range("D6").Value = "Flavia " + $nbsp
$ws.Range("E6").Value = "9.96,8.49"
$ws.Range("F6").Value = "Cyrus(mother): 0522363358"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "30.0"

# --- Row 7 (child 1) ---
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "16"
$ws.Range("C7").Value = "Collette " + $nbsp
$ws.Range("D7").Value = "Billi " + $nbsp
$ws.Range("E7").Value = "9.82,6.59"
$ws.Range("F7").Value = "Elias(mother): 0578741979"
$ws.Range("G7").Value = "7:03:00"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "27.0"

# --- Row 8 (child 2) ---
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "17"
$ws.Range("C8").Value = "Britta " + $nbsp
$ws.Range("D8").Value = "Jamel " + $nbsp
$ws.Range("E8").Value = "5.79,5.55"
$ws.Range("F8").Value = "Albertine(father): 0574981040"
$ws.Range("G8").Value = "7:08:00"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "22.0"

# --- Row 9 (child 3) ---
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "14"
$ws.Range("C9").Value = "Lorinda " + $nbsp
$ws.Range("D9").Value = "Tyron " + $nbsp
$ws.Range("E9").Value = "8.07,2.26"
$ws.Range("F9").Value = "Teresa(grandmother): 0558587699"
$ws.Range("G9").Value = "7:13:00"
# H9 stays "17.0" (unchanged)

# --- Row 10 (child 4) ---
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "1"
$ws.Range("C10").Value = "Corene " + $nbsp
$ws.Range("D10").Value = "Myra " + $nbsp
$ws.Range("E10").Value = "8.46,-1.79"
$ws.Range("F10").Value = "Georgie(mother): 0544823581"
$ws.Range("G10").Value = "7:19:00"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "11.0"

# --- Row 11 (school row) ---
$ws.Range("G11").Value = "7:30:00"
$ws.Range("H11").ClearContents()

# --- Row 13 (time row) ---
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "30.0"
